# Apply edits described by the commit: "metrics and summary of training"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text tweaks ---
$ws.Range("K2").Value = "Final Accuracy(%)"
$ws.Range("A4").Value = "V1.0"

# --- Expand the "Final Accuracy" merged header from K2:L2 to K2:M2,
#     add a Test column under it (M3), and move "Remarks" from M2 to N2 ---
$ws.Range("K2:L2").UnMerge()
$ws.Range("M2").Value = $ws.Range("M2").Value
$ws.Range("N2").Value = $ws.Range("M2").Value
$ws.Range("K2:M2").Merge()
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = "Remarks"
$ws.Range("M3").Value = "Test"

# --- Row 4 now only holds the run identification columns; its metrics
#     move to row 5 together with a duplicated set of descriptive columns. ---
$ws.Range("G4").Clear()
$ws.Range("H4").Clear()
$ws.Range("I4").Clear()
$ws.Range("J4").Clear()
$ws.Range("K4").Clear()
$ws.Range("L4").Clear()

# --- New row 5: V1.0 metrics ---
$ws.Range("A5").Value = "V1.1"
$ws.Range("B5").Value = 10
$ws.Range("C5").Value = "ResNet50"
$ws.Range("D5").Value = 16
$ws.Range("E5").Value = "Adam"
$ws.Range("F5").Value = "Cross Entropy "
$ws.Range("G5").Value = 0.000015
$ws.Range("H5").Value = 58647
$ws.Range("I5").Value = 19550
$ws.Range("J5").Value = 27599
$ws.Range("K5").Value = 99.88
$ws.Range("L5").Value = 98.93
$ws.Range("M5").Value = 98.93

# --- New row 6: V1.2 (15 epochs, lower LR) ---
$ws.Range("A6").Value = "V1.2"
$ws.Range("B6").Value = 15
$ws.Range("C6").Value = "ResNet50"
$ws.Range("D6").Value = 16
$ws.Range("E6").Value = "Adam"
$ws.Range("F6").Value = "Cross Entropy "
$ws.Range("G6").Value = 0.0000075
$ws.Range("H6").Value = 58647
$ws.Range("I6").Value = 19550
$ws.Range("J6").Value = 27599
$ws.Range("K6").Value = 98.64
$ws.Range("L6").Value = 97.6
$ws.Range("M6").Value = 97.72
$ws.Range("N6").Value = "Learning rate is decreased so is accuracy"

# --- New row 7: V1.3 (10 epochs, higher LR) ---
$ws.Range("A7").Value = "V1.3"
$ws.Range("B7").Value = 10
$ws.Range("C7").Value = "ResNet50"
$ws.Range("D7").Value = 16
$ws.Range("E7").Value = "Adam"
$ws.Range("F7").Value = "Cross Entropy "
$ws.Range("G7").Value = 0.0000155
$ws.Range("H7").Value = 58647
$ws.Range("I7").Value = 19550
$ws.Range("J7").Value = 27599
$ws.Range("K7").Value = 99.91
$ws.Range("L7").Value = 98.94
$ws.Range("M7").Value = 99.06
$ws.Range("N7").Value = "Increased learning rate which also increased accuracy slightly."

# --- Extend title merge to new last column ---
$ws.Range("A1:M1").UnMerge()
$ws.Range("A1:N1").Merge()
